$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("A1").Value = "firstname"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "id"
$ws.Range("D1").Value = "company"
$ws.Range("E1").Value = "ort"

# --- Row 2 (Hans Muster) ---
$ws.Range("A2").Value = "Hans"
$ws.Range("B2").Value = "Muster"
$ws.Range("C2").Value = "001"
$ws.Range("E2").Value = "Zurich"

# --- Row 3 (Simon Gadient) ---
$ws.Range("A3").Value = "Simon"
$ws.Range("B3").Value = "Gadient"
$ws.Range("C3").Value = "002"
$ws.Range("D3").Value = "Web Essentials"
$ws.Range("E3").Value = "Phnom Penh"

# --- Selection matches the author's final cursor position ---
$ws.Range("B1").Select()
